# Update "想去人数" (column F) figures across all four sheets to match the
# newly generated data snapshot (gh-pages output at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 615
$ws.Range("F5").Value = 2712
$ws.Range("F7").Value = 206
$ws.Range("F9").Value = 264
$ws.Range("F10").Value = 6055
$ws.Range("F11").Value = 3
$ws.Range("F14").Value = 4943
$ws.Range("F16").Value = 94
$ws.Range("F17").Value = 10
$ws.Range("F18").Value = 2548
$ws.Range("F19").Value = 1323
$ws.Range("F20").Value = 1495
$ws.Range("F21").Value = 1206
$ws.Range("F22").Value = 280
$ws.Range("F23").Value = 112
$ws.Range("F24").Value = 124
$ws.Range("F25").Value = 1010
$ws.Range("F26").Value = 222
$ws.Range("F27").Value = 385
$ws.Range("F28").Value = 6
$ws.Range("F29").Value = 1343
$ws.Range("F30").Value = 7
$ws.Range("F31").Value = 2075
$ws.Range("F32").Value = 283
$ws.Range("F34").Value = 60
$ws.Range("F35").Value = 236
$ws.Range("F36").Value = 1463
$ws.Range("F38").Value = 1007
$ws.Range("F39").Value = 111
$ws.Range("F41").Value = 260
$ws.Range("F42").Value = 1716
$ws.Range("F43").Value = 2511
$ws.Range("F44").Value = 54
$ws.Range("F45").Value = 107
$ws.Range("F46").Value = 261

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 313
$ws.Range("F9").Value = 167
$ws.Range("F11").Value = 198
$ws.Range("F16").Value = 42
$ws.Range("F23").Value = 329
$ws.Range("F24").Value = 27
$ws.Range("F26").Value = 13

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 1440
$ws.Range("F9").Value = 1805
$ws.Range("F10").Value = 2397
$ws.Range("F11").Value = 798
$ws.Range("F12").Value = 699

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 615
$ws.Range("F7").Value = 2712
$ws.Range("F8").Value = 206
$ws.Range("F9").Value = 1440
$ws.Range("F10").Value = 2397
$ws.Range("F11").Value = 6055
$ws.Range("F12").Value = 798
$ws.Range("F16").Value = 4943
$ws.Range("F17").Value = 94
$ws.Range("F18").Value = 2548
$ws.Range("F19").Value = 1323
$ws.Range("F20").Value = 1495
$ws.Range("F21").Value = 1206
$ws.Range("F22").Value = 280
$ws.Range("F23").Value = 112
$ws.Range("F24").Value = 124
$ws.Range("F25").Value = 167
$ws.Range("F26").Value = 222
$ws.Range("F27").Value = 385
$ws.Range("F28").Value = 1343
$ws.Range("F29").Value = 2075
$ws.Range("F30").Value = 283
$ws.Range("F32").Value = 236
$ws.Range("F33").Value = 42
$ws.Range("F34").Value = 1463
$ws.Range("F36").Value = 111
$ws.Range("F40").Value = 260
$ws.Range("F41").Value = 27
$ws.Range("F42").Value = 1716
$ws.Range("F43").Value = 2511
$ws.Range("F44").Value = 107
$ws.Range("F45").Value = 261
